$d = $word.ActiveDocument

# --- Paragraph 1: title block (date line + paper-title line, separated by <w:br/>) ---
$ok = $d.Content.Find.Execute("⚡️🚀המאמר היומי של מייק -02.11.24: ⚡️🚀", $true, $false, $false, $false, $false, $true, 1, $false, "⚡️🚀המאמר היומי של מייק -01.11.24: ⚡️🚀", 2)
if (-not $ok) { Write-Output "WARNING: replace failed for title date" }
$ok = $d.Content.Find.Execute("Learning to Compress: Local Rank and Information Compression in Deep Neural Networks", $true, $false, $false, $false, $false, $true, 1, $false, "LLMs Are In-Context Reinforcement Learners", 2)
if (-not $ok) { Write-Output "WARNING: replace failed for paper title" }

# --- Paragraph 2 ---
$ok = $d.Content.Find.Execute("היום סוקרים מאמר כחול לבן למחצה (אחד המחברים משניים הוא ישראלי רביד שוורץ זיו) והם חוקרים נושא שמעניין אותי מאוד באופן אישי. הנושא הוא דחיסה של דאטה באמצעות רשתות נוירונים והוא גם מאוד קשור לעבודות של נפתלי תשבי האגדי בנושא צוואר בקבוק מידעי (information bottleneck או IB) וגם השערת יריעה (manifold hypothesis או MH) בנוגע לרשתות נוירונים עמוקות.", $true, $false, $false, $false, $false, $true, 1, $false, "אני אוהב מאמרים שמשלבים כמה שיטות של ML. אסקור היום אחד כזה המציע לשדך למידת in-context עם למידה באמצעות חיזוקים או בקצרה RL. למידת in-context היא יכולת של מודל שפה ללמוד משהו חדש מכמה דוגמאות בפרופמט ללא צורך בפיין טיון. יש לא מעט הסברים ליכולת די מפתיע זו ולפעמים יכולת זו נקראה emergent capabilities.", 2)
if (-not $ok) { Write-Output "WARNING: replace failed for paragraph 2" }

# --- Paragraph 3 ---
$ok = $d.Content.Find.Execute("MH טוענת שדאטה מהעולם האמיתי (כגון תמונות או טקסט) אינם מפוזרים באופן אחיד במרחב בעל מימד גבוה, אלא שוכנים על יריעה בעל מימד נמוך יותר. רשתות נוירונים עמוקות מצליחות היטב עם הדאטה הז כי הן לומדות לזהות ולנצל את המבנה של אותה יריעה, מה שמאפשר להן לבצע הכללה טובה למרות המורכבות העצומה של המרחב המקורי.", $true, $false, $false, $false, $false, $true, 1, $false, "עכשיו נשאלת השאלה: איך נוכל לבחור דוגמאת להדגמה שאנו מראים למודל שפה בפרומפט למקסום ביצועיי המודל? השאלה הזו לא מאוד טריויאלית ואין עליה כרגע תשובה חד משמעית. המחברים מציעים לגשת לבעיה זו דרך למידה עם חיזוקים (סוג של). השיטה הנאיבית היא פשוט לצבור דוגמאות עד שנגמר לנו את אורך חלון ההקשר של המודל. לכל דוגמא בהדגמה אנו שומרים בבאפר את השלישיה המכילה את הדוגמא (שאלה עצמה)ֿ, תשובת המודל ומשערך של איכות התשובה (או פשוט האם התשובה נכונה או לא). ואז באינפרנס פשוט לוקחים את הדוגמאות האלו בתור פרומפט.", 2)
if (-not $ok) { Write-Output "WARNING: replace failed for paragraph 3" }

# --- Paragraph 4 ---
# (direct Range.Text assignment: the replacement text contains a straight
#  apostrophe in "באץ'" that Find/Replace would auto-correct into a curly quote)
$d.Paragraphs.Item(4).Range.Text = "לטענת המחברים הגישה הנאיבית הזו לא עובדת משתי סיבות עיקריות. קודם כל שילוב מתמשך של אותם הפרומפטים לדוגמאות שונות מוביל לשונות גדולה בפלט של LLM (לפי המחקרים הקודמים עלולה להוביל לביצועים ירודים). הסיבה השניה טמונה בכך ששלישיות (שאלה, תשובה, לא נכון) מסבכות את המודל ולא מספקות לו מספיק מידע על איך היה צריך לענות נכון (ד״א בלמידה ניגודית יש בעיה דומה המצריכה כמות מאוד גדולה של דוגמאות שליליות בכל באץ' - כתבתי על זה לא מעט בסקירותיי)."

# --- Paragraph 5 ---
$ok = $d.Content.Find.Execute("נציין שמטריצות עם דרגה לא מלאה מהוות מרחב בעל מידה אפס במרחב של כל המטריצות (כמו הסתברות של כל מספר עם דוגמים יוניפורמית בין 0 ל 1). עקב כך המאמר מגדיר robust local rank או RLR שזה מספר ערכים סינגולריים (הכללה של הערכים העצמיים) של היעקוביאן שהם גדולים ממספר קטן אפסילון אך חיובי (נזכור עבור דרגה אמיתית צריך להחליף אפסילון ב 0). ", $true, $false, $false, $false, $false, $true, 1, $false, "עקב כך המחברים הציעו להכניס קצת ״אקראיות״ לבניית הפרומפטים (המחברים קוראים לזה אפיזודה בהתאם לטרמינולוגיה של RL - כל אפיזודה מורכבת מכמה שלישיות של שאלה, תשובה, נכונות התשובה) וגם להשתמש באפיזודות שקיבלו ציון ״נכון״. לכל דוגמא הם הציע קודם לדגום באקראי מהבאפר של אפיזודות בצורה אקראית ולהשתמש לכל דוגמא במדגם שונה של אפיזודות. כאמור שומרים רק את האפיזודות שבהם המודל צדק. כך פרומפט לכל שאילתה הופך להיות לא קבוע ומכיל רק דוגמאות עם תשובות נכונות. זה נקרא Explorative ICRL במאמר.", 2)
if (-not $ok) { Write-Output "WARNING: replace failed for paragraph 5" }

# --- Paragraph 6 ---
$ok = $d.Content.Find.Execute("אוקיי, מקווה ששרדתם את זה אז עכשיו מגיעים שני המשפטים העיקריים של המאמר. הם טוענים שברשתות עמוקות (מספר שכבות גבוה) בבעיות סיווג תמיד יהיה שכבה l שה-RLR יהיה נמוך מ-(פרופורציונלי לאפסילון בחזקה מינוס 2 ובנורמת אופרטור של מטריצת השכבה l (נורמת אופרטור זה הערך הסינגולרי הגבוה ביותר). הכוונה כאן לרשת שעושה התאמה מושלמת לדאטה האימון (עם מרג'ין 1 כלומר מצליחה להפריד בין הקטגוריות השונות בבטחה). משמעות המשפט היא שהרשת המאומנת דוחסת את הדאטה בשכבה l באופן אפקטיבי.", $true, $false, $false, $false, $false, $true, 1, $false, "כמובן ש Explorative ICRL לא יעיל חישובית כי כל פעם צריך לחשב את הפרומפט מחדש (מה שלא צריך לעשות בגישה הנאיבית אך לא עובדת). המחברים שכללו את זה עם מנגנון קאשינג המאפשר לשמור מספר קבוע של פרומפרטים (מערך של אפיזודות) ולכל אפיזודה נתונה להחליט לאלו מהם להוסיף אותה. זה מקל על העלות החישובית.", 2)
if (-not $ok) { Write-Output "WARNING: replace failed for paragraph 6" }

# --- Paragraph 7 ---
$ok = $d.Content.Find.Execute("המחברים מוכיחים משפט דומה בנוגע לבעיות רגרסיה.", $true, $false, $false, $false, $false, $true, 1, $false, "מאמר חמוד למרות שמשום מה לקח לי קצת זמן להבין אותו…", 2)
if (-not $ok) { Write-Output "WARNING: replace failed for paragraph 7" }

# --- Paragraph 8 ("תמיד כיף...") is removed entirely ---
$p8 = $d.Paragraphs.Item(8)
if ($p8.Range.Text -like "*תמיד כיף*") { $p8.Range.Delete() } else { Write-Output "WARNING: paragraph 8 was not the expected one" }

# --- Former last paragraph (arxiv link) updated ---
$ok = $d.Content.Find.Execute("https://arxiv.org/abs/2410.07687", $true, $false, $false, $false, $false, $true, 1, $false, "https://arxiv.org/pdf/2410.05362", 2)
if (-not $ok) { Write-Output "WARNING: replace failed for arxiv link" }

